# crypto_analysis.xlsx update
#   - Sheet "Top 5 by Market Cap": insert a "symbol" column, refresh market
#     cap figures, append a "Last updated" footer row.
#   - Sheet "Average Price": refresh the average, append a "Last updated"
#     note next to it.
#   - Sheet "Highest Change" / "Lowest Change": rename header to "change",
#     replace the leader name + percent value with a "Last updated" note.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Top 5 by Market Cap"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a new column B ("symbol") before the existing market_cap column.
$ws1.Columns.Item(2).Insert()

$ws1.Range("B1").Value = "symbol"

$ws1.Range("B2").Value = "BTC"
$ws1.Range("C2").Value = 1872066201124.301

$ws1.Range("B3").Value = "ETH"
$ws1.Range("C3").Value = 398474056364.7794

$ws1.Range("B4").Value = "USDT"
$ws1.Range("C4").Value = 137556493146.6945

$ws1.Range("B5").Value = "XRP"
$ws1.Range("C5").Value = 133204662185.7822

$ws1.Range("B6").Value = "BNB"
$ws1.Range("C6").Value = 99835401046.23409

$ws1.Range("A7").Value = "Last updated (GMT): 2025-01-09 00:07:13"

# ---------------------------------------------------------------------------
# Sheet 2: "Average Price"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = 2009.267132097337
$ws2.Range("B2").Value = "Last updated (GMT): 2025-01-09 00:07:13"

# ---------------------------------------------------------------------------
# Sheet 3: "Highest Change"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("B1").Value = "change"
$ws3.Range("A2").Value = "MANTRA"
$ws3.Range("B2").Value = "Last updated (GMT): 2025-01-09 00:07:13"

# ---------------------------------------------------------------------------
# Sheet 4: "Lowest Change"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("B1").Value = "change"
$ws4.Range("A2").Value = "Ethena"
$ws4.Range("B2").Value = "Last updated (GMT): 2025-01-09 00:07:13"
